# Update "想去人数" (interested count) values in F column on two sheets:
# "展览" (Exhibitions) and "全部类型" (All types), matching the gh-pages
# data refresh (456a3b4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1132
$ws1.Range("F6").Value = 12114
$ws1.Range("F9").Value = 11888
$ws1.Range("F10").Value = 4776
$ws1.Range("F11").Value = 590
$ws1.Range("F12").Value = 84
$ws1.Range("F17").Value = 356

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1132
$ws4.Range("F8").Value = 12114
$ws4.Range("F11").Value = 11888
$ws4.Range("F12").Value = 4776
$ws4.Range("F13").Value = 590
$ws4.Range("F14").Value = 84
$ws4.Range("F19").Value = 356
